$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/date formatting from column F into the two new columns D and E
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (most recent two quarters) with their reported figures
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("D8").Value = 162000
$ws.Range("E8").Value = 161700
$ws.Range("D9").Value = 107900
$ws.Range("E9").Value = 105300
$ws.Range("D10").Value = 54100
$ws.Range("E10").Value = 56400
$ws.Range("D12").Value = 3300
$ws.Range("E12").Value = 3300
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 10400
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 164000
$ws.Range("E17").Value = 150400
$ws.Range("D18").Value = -2000
$ws.Range("E18").Value = 11300
$ws.Range("D20").Value = -2500
$ws.Range("E20").Value = -2300
$ws.Range("D21").Value = 9100
$ws.Range("E21").Value = 22200
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = -4500
$ws.Range("E23").Value = 9100
$ws.Range("D24").Value = 1200
$ws.Range("E24").Value = 2400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -5700
$ws.Range("E26").Value = 6700
$ws.Range("D27").Value = -5700
$ws.Range("E27").Value = 6700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2500
$ws.Range("E32").Value = 2300
$ws.Range("D33").Value = -5700
$ws.Range("E33").Value = 6700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -5700
$ws.Range("E35").Value = 6700
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("D41").Value = 37500
$ws.Range("E41").Value = 36400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 78800
$ws.Range("E43").Value = 65000
$ws.Range("D44").Value = 173500
$ws.Range("E44").Value = 175200
$ws.Range("D45").Value = 7600
$ws.Range("E45").Value = 8200
$ws.Range("D46").Value = 297400
$ws.Range("E46").Value = 284900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 185600
$ws.Range("E48").Value = 179600
$ws.Range("D49").Value = 279300
$ws.Range("E49").Value = 293800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 10500
$ws.Range("E52").Value = 11300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 772800
$ws.Range("E54").Value = 769600
$ws.Range("D57").Value = 33900
$ws.Range("E57").Value = 38400
$ws.Range("D58").Value = 6300
$ws.Range("E58").Value = 6300
$ws.Range("D59").Value = 61100
$ws.Range("E59").Value = 54300
$ws.Range("D60").Value = 101300
$ws.Range("E60").Value = 99000
$ws.Range("D61").Value = 221500
$ws.Range("E61").Value = 216000
$ws.Range("D62").Value = 18100
$ws.Range("E62").Value = 18300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 340900
$ws.Range("E66").Value = 333400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 393100
$ws.Range("E72").Value = 398800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 432000
$ws.Range("E76").Value = 436200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("D81").Value = -5700
$ws.Range("E81").Value = 6700
$ws.Range("D83").Value = 13600
$ws.Range("E83").Value = 13100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 11700
$ws.Range("E89").Value = -1600
$ws.Range("D91").Value = -6400
$ws.Range("E91").Value = -12700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -8300
$ws.Range("E94").Value = -11500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2300
$ws.Range("E100").Value = 24200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 1100
$ws.Range("E102").Value = 11200

# A handful of prior-quarter figures were also restated in this update
$ws.Range("F24").Value = 12500
$ws.Range("F26").Value = -1100
$ws.Range("F27").Value = -1100
$ws.Range("F29").Value = 8700
$ws.Range("F91").Value = -6900
$ws.Range("G91").Value = -4500
$ws.Range("H91").Value = -4100
$ws.Range("I91").Value = -5200
$ws.Range("J91").Value = -4700
